# Apply manufacturer / orderer list cleanup to testdata workbook.
$wb = $excel.ActiveWorkbook

$wsManufacturer = $wb.Worksheets.Item("Manufacturer")
$wsOrderer = $wb.Worksheets.Item("Orderer")

# ---------------------------------------------------------------------------
# 1. Remove a handful of manufacturer rows entirely (duplicates / stray
#    entries that were cleaned up from the master list).
# ---------------------------------------------------------------------------
$namesToRemove = @(
    "Fjällräven",
    "BOND NO. 9",
    "C.P. Company",
    "J.Crew",
    "Kiehl's",
    "Maison Kitsuné",
    "M·A·C",
    "NN07",
    "Samsøe & Samsøe"
)

foreach ($name in $namesToRemove) {
    $found = $wsManufacturer.Columns.Item(1).Find($name)
    if ($found -ne $null) {
        $found.EntireRow.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. Normalise a number of manufacturer names that contained accented /
#    special characters (apostrophes, diacritics, interpuncts) to their
#    plain-ASCII equivalents.
# ---------------------------------------------------------------------------
$manufacturerRenames = @{
    "Bric's"         = "Brics"
    "Chasin'"        = "Chasin"
    "Church's"       = "Churchs"
    "Claesen's"      = "Claesens"
    "Drake's"        = "Drake"
    "Estée Lauder"   = "Estee Lauder"
    "HERMÈS"         = "HERMES"
    "Lancôme"        = "Lancome"
    "Tumble 'n Dry"  = "Tumble n Dry"
    "Wood'd"         = "Wood"
}

foreach ($old in $manufacturerRenames.Keys) {
    $new = $manufacturerRenames[$old]
    $found = $wsManufacturer.Columns.Item(1).Find($old)
    if ($found -ne $null) {
        $found.Value = $new
    }
}

# Same normalisation for the Orderer sheet.
$ordererRenames = @{
    "Hunkemöller" = "Hunkemoller"
}

foreach ($old in $ordererRenames.Keys) {
    $new = $ordererRenames[$old]
    $found = $wsOrderer.Columns.Item(1).Find($old)
    if ($found -ne $null) {
        $found.Value = $new
    }
}

# ---------------------------------------------------------------------------
# 3. Re-sort both lists alphabetically (case-insensitive), keeping the
#    header row in place.
# ---------------------------------------------------------------------------
$manufacturerLastRow = $wsManufacturer.Cells.Item($wsManufacturer.Rows.Count, 1).End(-4162).Row
$manufacturerRange = $wsManufacturer.Range("A1:A" + $manufacturerLastRow)
$manufacturerRange.Sort($wsManufacturer.Range("A2"), 1, $null, $null, 1, $null, 1, 1)

$ordererLastRow = $wsOrderer.Cells.Item($wsOrderer.Rows.Count, 1).End(-4162).Row
$ordererRange = $wsOrderer.Range("A1:A" + $ordererLastRow)
$ordererRange.Sort($wsOrderer.Range("A2"), 1, $null, $null, 1, $null, 1, 1)

# ---------------------------------------------------------------------------
# 4. Restore the view state (scroll position / selection) Excel saved the
#    workbook with. The Orderer sheet is touched first so that the
#    Manufacturer sheet ends up as the active (tab-selected) sheet again.
# ---------------------------------------------------------------------------
$wsOrderer.Activate()
$excel.ActiveWindow.Zoom = 100
$wsOrderer.Range("A25").Select()

$wsManufacturer.Activate()
$excel.ActiveWindow.ScrollRow = 293
$wsManufacturer.Range("A3:XFD3").Select()
